# "add dev part to scrum" -- assign developers to the Sprint backlog items,
# tidy up the row heights now that the Developer column holds short names,
# set the print orientation for the sheet and leave the selection where the
# user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint backlog")

# Column G = "Developer" -- fill in the actual team member assigned to each item.
$ws.Range("G2").Value = "XIAO KAI"
$ws.Range("G3").Value = "Yuyi"
$ws.Range("G4").Value = "SUN Hao"
$ws.Range("G5").Value = "XIAO KAI"
$ws.Range("G6").Value = "Maitre Robin"
$ws.Range("G7").Value = "Yann Vaillant"
$ws.Range("G8").Value = "FENG Changhong"
$ws.Range("G9").Value = "Yuyi"
$ws.Range("G10").Value = "SUN Hao"

# The developer names are shorter than the previous placeholder text, so a
# few rows no longer need to be as tall.
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 45

# Set the sheet up for printing in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection on the cell the user last touched.
$ws.Range("L6").Select()
